$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.61"
$ws.Range("D3").Value = "'23.00"
$ws.Range("D4").Value = "'5.388"
$ws.Range("D6").Value = "'3.426"
$ws.Range("D7").Value = "'6.502"
$ws.Range("D8").Value = "'0.8128"
$ws.Range("D9").Value = "'0.9253"
$ws.Range("D10").Value = "'0.1427"
$ws.Range("D11").Value = "'0.07412"
$ws.Range("D12").Value = "'0.03277"
$ws.Range("D13").Value = "'0.03086"
$ws.Range("D14").Value = "'0.09348"
$ws.Range("D15").Value = "'3.855"
$ws.Range("D16").Value = "'0.001570"
$ws.Range("D18").Value = "'0.0005988"
$ws.Range("D19").Value = "'0.005874"
$ws.Range("D20").Value = "'0.001261"
$ws.Range("E20").Value = "19BitKanKANBestin24h"
$ws.Range("D21").Value = "'0.004789"
$ws.Range("D22").Value = "'0.00007998"
$ws.Range("D24").Value = "'2.133"
$ws.Range("D25").Value = "'0.3238"
$ws.Range("D27").Value = "'0.0002339"
$ws.Range("D40").Value = "'0.03933"
$ws.Range("D41").Value = "'0.006298"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1074"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002619"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.008894"
$ws.Range("D45").Value = "'0.00005171"
$ws.Range("D47").Value = "'0.6897"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D50").Value = "'0.0001999"
